# Update the "Förändrad" date column (column C, rows 2-6) from
# 2023-10-13 (serial 45212) to 2023-10-22 (serial 45221) on the
# "Avverkningsanmälningar" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

$ws.Range("C2:C6").Value = 45221
